# Auto-generated Excel COM-interop script
# Applies cached numeric value updates to Sheets per the target diff.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3149.6428  # H19 (3646.913 -> 3149.6428)
$ws.Cells.Item(19, 10).Value = 3075.7878  # J19 (3792.5 -> 3075.7878)
$ws.Cells.Item(19, 12).Value = 3075.7878  # L19 (3792.5 -> 3075.7878)
$ws.Cells.Item(19, 14).Value = -3425.7878  # N19 (-4142.5 -> -3425.7878)
$ws.Cells.Item(28, 8).Value = 1628.1818  # H28 (1851.7368 -> 1628.1818)
$ws.Cells.Item(28, 9).Value = 312.94116  # I28 (338.86667 -> 312.94116)
$ws.Cells.Item(28, 10).Value = 6100  # J28 (7525 -> 6100)
$ws.Cells.Item(28, 11).Value = 312.94116  # K28 (338.86667 -> 312.94116)
$ws.Cells.Item(28, 12).Value = 6100  # L28 (7525 -> 6100)
$ws.Cells.Item(28, 13).Value = 172.05884  # M28 (146.13333 -> 172.05884)
$ws.Cells.Item(28, 14).Value = -7070  # N28 (-8495 -> -7070)
$ws.Cells.Item(62, 8).Value = 4916.6665  # H62 (4690.4287 -> 4916.6665)
$ws.Cells.Item(62, 9).Value = 4666.6665  # I62 (4333.25 -> 4666.6665)
$ws.Cells.Item(62, 11).Value = 4666.6665  # K62 (4333.25 -> 4666.6665)
$ws.Cells.Item(62, 13).Value = -4042.6665  # M62 (-3709.25 -> -4042.6665)
$ws.Cells.Item(64, 8).Value = 6129.091  # H64 (5168.9165 -> 6129.091)
$ws.Cells.Item(64, 9).Value = 5908.6665  # I64 (5631.25 -> 5908.6665)
$ws.Cells.Item(64, 10).Value = 6211.75  # J64 (4937.75 -> 6211.75)
$ws.Cells.Item(64, 11).Value = 5908.6665  # K64 (5631.25 -> 5908.6665)
$ws.Cells.Item(64, 12).Value = 6211.75  # L64 (4937.75 -> 6211.75)
$ws.Cells.Item(64, 13).Value = -5660.6665  # M64 (-5383.25 -> -5660.6665)
$ws.Cells.Item(64, 14).Value = -6707.75  # N64 (-5433.75 -> -6707.75)
$ws.Cells.Item(65, 8).Value = 4916.6665  # H65 (4690.4287 -> 4916.6665)
$ws.Cells.Item(65, 9).Value = 4666.6665  # I65 (4333.25 -> 4666.6665)
$ws.Cells.Item(65, 11).Value = 23333.3325  # K65 (21666.25 -> 23333.3325)
$ws.Cells.Item(65, 13).Value = -20213.3325  # M65 (-18546.25 -> -20213.3325)
$ws.Cells.Item(67, 8).Value = 6129.091  # H67 (5168.9165 -> 6129.091)
$ws.Cells.Item(67, 9).Value = 5908.6665  # I67 (5631.25 -> 5908.6665)
$ws.Cells.Item(67, 10).Value = 6211.75  # J67 (4937.75 -> 6211.75)
$ws.Cells.Item(67, 11).Value = 5908.6665  # K67 (5631.25 -> 5908.6665)
$ws.Cells.Item(67, 12).Value = 6211.75  # L67 (4937.75 -> 6211.75)
$ws.Cells.Item(67, 13).Value = -5050.6665  # M67 (-4773.25 -> -5050.6665)
$ws.Cells.Item(67, 14).Value = -7927.75  # N67 (-6653.75 -> -7927.75)
$ws.Cells.Item(86, 8).Value = 4488.25  # H86 (4738.25 -> 4488.25)
$ws.Cells.Item(86, 9).Value = 4984.3335  # I86 (6476.5 -> 4984.3335)
$ws.Cells.Item(86, 11).Value = 4984.3335  # K86 (6476.5 -> 4984.3335)
$ws.Cells.Item(86, 13).Value = -3861.3335  # M86 (-5353.5 -> -3861.3335)
$ws.Cells.Item(89, 8).Value = 4488.25  # H89 (4738.25 -> 4488.25)
$ws.Cells.Item(89, 9).Value = 4984.3335  # I89 (6476.5 -> 4984.3335)
$ws.Cells.Item(89, 11).Value = 24921.6675  # K89 (32382.5 -> 24921.6675)
$ws.Cells.Item(89, 13).Value = -19305.6675  # M89 (-26766.5 -> -19305.6675)
$ws.Cells.Item(92, 8).Value = 1495.8125  # H92 (1577.7142 -> 1495.8125)
$ws.Cells.Item(92, 9).Value = 1358.1538  # I92 (1437.3636 -> 1358.1538)
$ws.Cells.Item(92, 11).Value = 1358.1538  # K92 (1437.3636 -> 1358.1538)
$ws.Cells.Item(92, 13).Value = -110.1538  # M92 (-189.3635999999999 -> -110.1538)
$ws.Cells.Item(96, 8).Value = 15471.143  # H96 (14043.5 -> 15471.143)
$ws.Cells.Item(96, 9).Value = 25324.75  # I96 (25337.25 -> 25324.75)
$ws.Cells.Item(96, 10).Value = 2333  # J96 (2749.75 -> 2333)
$ws.Cells.Item(96, 11).Value = 75974.25  # K96 (76011.75 -> 75974.25)
$ws.Cells.Item(96, 12).Value = 6999  # L96 (8249.25 -> 6999)
$ws.Cells.Item(96, 13).Value = -74601.25  # M96 (-74638.75 -> -74601.25)
$ws.Cells.Item(96, 14).Value = -9745  # N96 (-10995.25 -> -9745)
$ws.Cells.Item(98, 8).Value = 1771.3684  # H98 (1831.5555 -> 1771.3684)
$ws.Cells.Item(98, 9).Value = 1203.2222  # I98 (1233.5294 -> 1203.2222)
$ws.Cells.Item(98, 11).Value = 1203.2222  # K98 (1233.5294 -> 1203.2222)
$ws.Cells.Item(98, 13).Value = 294.7778000000001  # M98 (264.4706000000001 -> 294.7778000000001)
$ws.Cells.Item(103, 8).Value = 789.6316  # H103 (789.2 -> 789.6316)
$ws.Cells.Item(103, 9).Value = 781.625  # I103 (781.55554 -> 781.625)
$ws.Cells.Item(103, 11).Value = 2344.875  # K103 (2344.66662 -> 2344.875)
$ws.Cells.Item(103, 13).Value = -1758.875  # M103 (-1758.66662 -> -1758.875)
$ws.Cells.Item(107, 8).Value = 242  # H107 (246.28572 -> 242)
$ws.Cells.Item(107, 9).Value = 206.2  # I107 (219.07692 -> 206.2)
$ws.Cells.Item(107, 11).Value = 206.2  # K107 (219.07692 -> 206.2)
$ws.Cells.Item(107, 13).Value = 1713.8  # M107 (1700.92308 -> 1713.8)
$ws.Cells.Item(115, 8).Value = 2997.8333  # H115 (2611.2144 -> 2997.8333)
$ws.Cells.Item(115, 9).Value = 710.875  # I115 (627 -> 710.875)
$ws.Cells.Item(115, 11).Value = 2132.625  # K115 (1881 -> 2132.625)
$ws.Cells.Item(115, 13).Value = -565.625  # M115 (-314 -> -565.625)
$ws.Cells.Item(118, 8).Value = 1434  # H118 (1434.0625 -> 1434)
$ws.Cells.Item(118, 9).Value = 1442.4286  # I118 (1442.5 -> 1442.4286)
$ws.Cells.Item(118, 11).Value = 4327.2858  # K118 (4327.5 -> 4327.2858)
$ws.Cells.Item(118, 13).Value = -2670.2858  # M118 (-2670.5 -> -2670.2858)
$ws.Cells.Item(122, 8).Value = 1771.3684  # H122 (1831.5555 -> 1771.3684)
$ws.Cells.Item(122, 9).Value = 1203.2222  # I122 (1233.5294 -> 1203.2222)
$ws.Cells.Item(122, 11).Value = 3609.6666  # K122 (3700.5882 -> 3609.6666)
$ws.Cells.Item(122, 13).Value = -1159.6666  # M122 (-1250.5882 -> -1159.6666)
$ws.Cells.Item(137, 8).Value = 2667.5  # H137 (2562.9778 -> 2667.5)
$ws.Cells.Item(137, 9).Value = 2885.6667  # I137 (2997.1177 -> 2885.6667)
$ws.Cells.Item(137, 10).Value = 2503.875  # J137 (2299.3928 -> 2503.875)
$ws.Cells.Item(137, 11).Value = 8657.000100000001  # K137 (8991.3531 -> 8657.000100000001)
$ws.Cells.Item(137, 12).Value = 7511.625  # L137 (6898.178400000001 -> 7511.625)
$ws.Cells.Item(137, 13).Value = -6107.000100000001  # M137 (-6441.3531 -> -6107.000100000001)
$ws.Cells.Item(137, 14).Value = -12611.625  # N137 (-11998.1784 -> -12611.625)
$ws.Cells.Item(138, 8).Value = 2895.9539  # H138 (2938.6833 -> 2895.9539)
$ws.Cells.Item(138, 9).Value = 2266  # I138 (2272.0908 -> 2266)
$ws.Cells.Item(138, 10).Value = 3196.6135  # J138 (3324.6052 -> 3196.6135)
$ws.Cells.Item(138, 11).Value = 6798  # K138 (6816.2724 -> 6798)
$ws.Cells.Item(138, 12).Value = 9589.8405  # L138 (9973.8156 -> 9589.8405)
$ws.Cells.Item(138, 13).Value = -1658  # M138 (-1676.2724 -> -1658)
$ws.Cells.Item(138, 14).Value = -19869.8405  # N138 (-20253.8156 -> -19869.8405)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 10414.5  # H28 (10040.363 -> 10414.5)
$ws.Cells.Item(28, 9).Value = 11043.125  # I28 (11218 -> 11043.125)
$ws.Cells.Item(28, 10).Value = 7900  # J28 (6900 -> 7900)
$ws.Cells.Item(28, 11).Value = 11043.125  # K28 (11218 -> 11043.125)
$ws.Cells.Item(28, 12).Value = 7900  # L28 (6900 -> 7900)
$ws.Cells.Item(28, 13).Value = -10851.125  # M28 (-11026 -> -10851.125)
$ws.Cells.Item(28, 14).Value = -8284  # N28 (-7284 -> -8284)
$ws.Cells.Item(32, 8).Value = 13587.849  # H32 (14396.774 -> 13587.849)
$ws.Cells.Item(32, 9).Value = 7765.077  # I32 (8324.708000000001 -> 7765.077)
$ws.Cells.Item(32, 11).Value = 7765.077  # K32 (8324.708000000001 -> 7765.077)
$ws.Cells.Item(32, 13).Value = -7478.077  # M32 (-8037.708000000001 -> -7478.077)
$ws.Cells.Item(99, 8).Value = 10414.5  # H99 (10040.363 -> 10414.5)
$ws.Cells.Item(99, 9).Value = 11043.125  # I99 (11218 -> 11043.125)
$ws.Cells.Item(99, 10).Value = 7900  # J99 (6900 -> 7900)
$ws.Cells.Item(99, 11).Value = 11043.125  # K99 (11218 -> 11043.125)
$ws.Cells.Item(99, 12).Value = 7900  # L99 (6900 -> 7900)
$ws.Cells.Item(99, 13).Value = -8048.125  # M99 (-8223 -> -8048.125)
$ws.Cells.Item(99, 14).Value = -13890  # N99 (-12890 -> -13890)
$ws.Cells.Item(130, 8).Value = 39000  # H130 (0 -> 39000)
$ws.Cells.Item(130, 10).Value = 39000  # J130 (0 -> 39000)
$ws.Cells.Item(130, 12).Value = 39000  # L130 (0 -> 39000)
$ws.Cells.Item(130, 14).Value = -49040  # N130 (None -> -49040)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(130, 8).Value = 89545.45  # H130 (90000 -> 89545.45)
$ws.Cells.Item(130, 10).Value = 89545.45  # J130 (90000 -> 89545.45)
$ws.Cells.Item(130, 12).Value = 89545.45  # L130 (90000 -> 89545.45)
$ws.Cells.Item(130, 14).Value = -99585.45  # N130 (-100040 -> -99585.45)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(130, 8).Value = 85066.664  # H130 (84062.5 -> 85066.664)
$ws.Cells.Item(130, 10).Value = 85066.664  # J130 (84062.5 -> 85066.664)
$ws.Cells.Item(130, 12).Value = 85066.664  # L130 (84062.5 -> 85066.664)
$ws.Cells.Item(130, 14).Value = -95106.664  # N130 (-94102.5 -> -95106.664)
$ws.Cells.Item(139, 8).Value = 92299.88  # H139 (92437.375 -> 92299.88)
$ws.Cells.Item(139, 10).Value = 92299.88  # J139 (92437.375 -> 92299.88)
$ws.Cells.Item(139, 12).Value = 92299.88  # L139 (92437.375 -> 92299.88)
$ws.Cells.Item(139, 14).Value = -102579.88  # N139 (-102717.375 -> -102579.88)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1466.0625  # H14 (1750.5 -> 1466.0625)
$ws.Cells.Item(14, 9).Value = 1466.0625  # I14 (1750.5 -> 1466.0625)
$ws.Cells.Item(14, 11).Value = 4398.1875  # K14 (5251.5 -> 4398.1875)
$ws.Cells.Item(14, 13).Value = -4225.1875  # M14 (-5078.5 -> -4225.1875)
$ws.Cells.Item(69, 8).Value = 1065.3334  # H69 (1178.6 -> 1065.3334)
$ws.Cells.Item(69, 9).Value = 1065.3334  # I69 (1178.6 -> 1065.3334)
$ws.Cells.Item(69, 11).Value = 3196.0002  # K69 (3535.8 -> 3196.0002)
$ws.Cells.Item(69, 13).Value = -2385.0002  # M69 (-2724.8 -> -2385.0002)
$ws.Cells.Item(72, 8).Value = 1065.3334  # H72 (1178.6 -> 1065.3334)
$ws.Cells.Item(72, 9).Value = 1065.3334  # I72 (1178.6 -> 1065.3334)
$ws.Cells.Item(72, 11).Value = 9588.000599999999  # K72 (10607.4 -> 9588.000599999999)
$ws.Cells.Item(72, 13).Value = -5532.000599999999  # M72 (-6551.4 -> -5532.000599999999)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 54462.844  # H98 (57043.89 -> 54462.844)
$ws.Cells.Item(98, 10).Value = 54462.844  # J98 (57043.89 -> 54462.844)
$ws.Cells.Item(98, 12).Value = 54462.844  # L98 (57043.89 -> 54462.844)
$ws.Cells.Item(98, 14).Value = -60452.844  # N98 (-63033.89 -> -60452.844)
$ws.Cells.Item(102, 8).Value = 1974.6  # H102 (1917.5714 -> 1974.6)
$ws.Cells.Item(102, 10).Value = 3097.5  # J102 (2324 -> 3097.5)
$ws.Cells.Item(102, 12).Value = 3097.5  # L102 (2324 -> 3097.5)
$ws.Cells.Item(102, 14).Value = -6341.5  # N102 (-5568 -> -6341.5)
$ws.Cells.Item(126, 8).Value = 5433  # H126 (5674.75 -> 5433)
$ws.Cells.Item(126, 9).Value = 5149.6665  # I126 (5479.8 -> 5149.6665)
$ws.Cells.Item(126, 11).Value = 15448.9995  # K126 (16439.4 -> 15448.9995)
$ws.Cells.Item(126, 13).Value = -12978.9995  # M126 (-13969.4 -> -12978.9995)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3956.1428  # H7 (4182.5 -> 3956.1428)
$ws.Cells.Item(7, 10).Value = 2598  # J7 (0 -> 2598)
$ws.Cells.Item(7, 12).Value = 2598  # L7 (0 -> 2598)
$ws.Cells.Item(7, 14).Value = -2822  # N7 (None -> -2822)
$ws.Cells.Item(40, 8).Value = 6431.6665  # H40 (6805.4614 -> 6431.6665)
$ws.Cells.Item(40, 9).Value = 6469.643  # I40 (6880.9165 -> 6469.643)
$ws.Cells.Item(40, 11).Value = 6469.643  # K40 (6880.9165 -> 6469.643)
$ws.Cells.Item(40, 13).Value = -6333.643  # M40 (-6744.9165 -> -6333.643)
$ws.Cells.Item(61, 8).Value = 2231.8333  # H61 (2257.1667 -> 2231.8333)
$ws.Cells.Item(61, 9).Value = 2162  # I61 (2189.6365 -> 2162)
$ws.Cells.Item(61, 11).Value = 2162  # K61 (2189.6365 -> 2162)
$ws.Cells.Item(61, 13).Value = -1960  # M61 (-1987.6365 -> -1960)
$ws.Cells.Item(100, 8).Value = 3567.4614  # H100 (3782.4 -> 3567.4614)
$ws.Cells.Item(100, 9).Value = 3567.4614  # I100 (3782.4 -> 3567.4614)
$ws.Cells.Item(100, 11).Value = 3567.4614  # K100 (3782.4 -> 3567.4614)
$ws.Cells.Item(100, 13).Value = -3026.4614  # M100 (-3241.4 -> -3026.4614)
$ws.Cells.Item(113, 8).Value = 2231.8333  # H113 (2257.1667 -> 2231.8333)
$ws.Cells.Item(113, 9).Value = 2162  # I113 (2189.6365 -> 2162)
$ws.Cells.Item(113, 11).Value = 2162  # K113 (2189.6365 -> 2162)
$ws.Cells.Item(113, 13).Value = 8  # M113 (-19.63650000000007 -> 8)
$ws.Cells.Item(126, 8).Value = 3956.1428  # H126 (4182.5 -> 3956.1428)
$ws.Cells.Item(126, 10).Value = 2598  # J126 (0 -> 2598)
$ws.Cells.Item(126, 12).Value = 7794  # L126 (0 -> 7794)
$ws.Cells.Item(126, 14).Value = -12734  # N126 (None -> -12734)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 412.33334  # H107 (365.8 -> 412.33334)
$ws.Cells.Item(107, 9).Value = 358.3125  # I107 (310.95 -> 358.3125)
$ws.Cells.Item(107, 11).Value = 1074.9375  # K107 (932.8499999999999 -> 1074.9375)
$ws.Cells.Item(107, 13).Value = 845.0625  # M107 (987.1500000000001 -> 845.0625)
$ws.Cells.Item(113, 8).Value = 479  # H113 (601.6667 -> 479)
$ws.Cells.Item(113, 9).Value = 479  # I113 (692.4 -> 479)
$ws.Cells.Item(113, 10).Value = 0  # J113 (148 -> 0)
$ws.Cells.Item(113, 11).Value = 1437  # K113 (2077.2 -> 1437)
$ws.Cells.Item(113, 12).Value = 0  # L113 (444 -> 0)
$ws.Cells.Item(113, 13).Value = 733  # M113 (92.80000000000018 -> 733)
$ws.Cells.Item(113, 14).ClearContents()  # N113 (was -4784)
$ws.Cells.Item(132, 8).Value = 3920.9143  # H132 (4111.0312 -> 3920.9143)
$ws.Cells.Item(132, 9).Value = 4108.1724  # I132 (4363.769 -> 4108.1724)
$ws.Cells.Item(132, 11).Value = 12324.5172  # K132 (13091.307 -> 12324.5172)
$ws.Cells.Item(132, 13).Value = -9794.517200000002  # M132 (-10561.307 -> -9794.517200000002)
$ws.Cells.Item(136, 8).Value = 28125.182  # H136 (29403.047 -> 28125.182)
$ws.Cells.Item(136, 9).Value = 30187.7  # I136 (31708.63 -> 30187.7)
$ws.Cells.Item(136, 11).Value = 90563.10000000001  # K136 (95125.89 -> 90563.10000000001)
$ws.Cells.Item(136, 13).Value = -88013.10000000001  # M136 (-92575.89 -> -88013.10000000001)
